$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spell-check / typography fixes across the "Activité" (and resume) column.

$ws.Range("B5").Value  = "Introduction du module. Constitution du groupe avec nomination d’un chef de groupe et de son remplaçant."
$ws.Range("B6").Value  = "Réunion: choix de proposition projet"
$ws.Range("B7").Value  = "Feedback des propositions. Discussion pour compléter le projet choisi: DARYLL"
$ws.Range("B8").Value  = "Rédaction du modèle de Gantt sous Excel, ainsi que des répartitions des heures"
$ws.Range("B9").Value  = "Réunion: Discuter des fonctionnalités du projet"
$ws.Range("B11").Value = "Changement de modèle de Gantt"
$ws.Range("B13").Value = "Réunion: Refaire la planification initiale"
$ws.Range("B14").Value = "Réunion: Définir les tâches à réaliser des vacances"
$ws.Range("B16").Value = "Réunion: Définir les tâches à réaliser des vacances"
$ws.Range("B19").Value = "Création des classes pour la communication Java — MySQL"
$ws.Range("B20").Value = "Réunion: avancement après vacances"
$ws.Range("B27").Value = "Réunion de groupe afin de définir les tâches à faire durant la période de l’ascension"
$ws.Range("B30").Value = "Réunion du petit groupe: Restructuration du projet"
$ws.Range("B32").Value = "Correction de bug Java sur les calls des méthodes SQL"

# Added the resume (name) for Aurélien Manalito Adrien Siu.
$ws.Range("A2").Value = "Dejvid Muaremi"

# Move the active selection to E22, as recorded in the saved view state.
$ws.Range("E22").Select()
